# coordinadoresTorresRegiones.xlsx
# "agregado PabloSanMartin, actualizado gitignore"
#
# MARTIN PALMA (id YP11856) is replaced by PABLO SEBASTIAN SAN MARTIN
# (id SE45933) as coordinator for REGION=CENTRONORTE, in both the
# COMUNICACIONES row (row 2) and the INFRAESTRUCTURA row (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "PABLO SEBASTIAN SAN MARTIN"
$ws.Range("D2").Value = "SE45933"
$ws.Range("C3").Value = "PABLO SEBASTIAN SAN MARTIN"
$ws.Range("D3").Value = "SE45933"

# Leave the selection where the editor's session ended up.
$ws.Range("H17").Select() | Out-Null
